# Refresh the cryptocurrency snapshot table on Sheet1 (Coin | Link | Price | Volume(1h))
# with the latest scraped price/volume figures. Cells D/E hold formatted text
# (prices such as "1.01" and percentages such as "  +0.53%  ") rather than numbers,
# so numeric-looking values are written with a leading apostrophe to force text,
# then ClearFormats() strips the resulting 'text number format' styling so the
# cell ends up as a plain, unstyled text cell (matching the rest of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.980.05"
$ws.Range("E2").Value = "  +0.53%  "
# Row 3
$ws.Range("D3").Value = "1.557.09"
$ws.Range("E3").Value = "  -0.18%  "
# Row 4
$ws.Range("D4").Value = "'1.01"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.41%  "
# Row 5
$ws.Range("D5").Value = "'206.98"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.76%  "
# Row 6
$ws.Range("E6").Value = "  +1.11%  "
# Row 7
$ws.Range("E7").Value = "  +0.40%  "
# Row 8
$ws.Range("E8").Value = "  +0.77%  "
# Row 9
$ws.Range("E9").Value = "  -0.15%  "
# Row 10
$ws.Range("E10").Value = "  -0.22%  "
# Row 11
$ws.Range("D11").Value = "'0.0858"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.40%  "
# Row 12
$ws.Range("D12").Value = "1.778.80"
$ws.Range("E12").Value = "  -0.15%  "
# Row 13
$ws.Range("D13").Value = "1.550.87"
$ws.Range("E13").Value = "  -0.92%  "
# Row 14
$ws.Range("E14").Value = "  -0.38%  "
# Row 15
$ws.Range("E15").Value = "  +0.36%  "
# Row 16
$ws.Range("D16").Value = "26.966.31"
$ws.Range("E16").Value = "  +0.50%  "
# Row 17
$ws.Range("D17").Value = "'61.78"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.87%  "
# Row 18
$ws.Range("D18").Value = "'214.88"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.09%  "
# Row 19
$ws.Range("D19").Value = "0.0₃0686"
$ws.Range("E19").Value = "  +0.80%  "
# Row 20
$ws.Range("D20").Value = "'7.25"
$ws.Range("D20").ClearFormats()
# Row 21
$ws.Range("D21").Value = "'1.01"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.39%  "
# Row 22
$ws.Range("E22").Value = "  -1.77%  "
# Row 23
$ws.Range("D23").Value = "'9.22"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.75%  "
# Row 24
$ws.Range("D24").Value = "'1.97"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.15%  "
# Row 25
$ws.Range("D25").Value = "'153.46"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.10%  "
# Row 26
$ws.Range("D26").Value = "'6.64"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.19%  "
# Row 27
$ws.Range("E27").Value = "  -0.80%  "
# Row 28
$ws.Range("E28").Value = "  +0.38%  "
# Row 29
$ws.Range("E29").Value = "  +0.57%  "
# Row 30
$ws.Range("E30").Value = "  -1.15%  "
# Row 31
$ws.Range("E31").Value = "  -0.48%  "
# Row 32
$ws.Range("E32").Value = "  +2.07%  "
# Row 33
$ws.Range("D33").Value = "1.373.90"
$ws.Range("E33").Value = "  -0.55%  "
# Row 34
$ws.Range("E34").Value = "  +1.40%  "
# Row 35
$ws.Range("D35").Value = "'1.55"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.83%  "
# Row 36
$ws.Range("D36").Value = "'0.968"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.56%  "
# Row 37
$ws.Range("E37").Value = "  +0.61%  "
# Row 38
$ws.Range("E38").Value = "  +0.88%  "
# Row 39
$ws.Range("D39").Value = "'0.521"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.24%  "
# Row 40
$ws.Range("D40").Value = "'0.809"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.05%  "
# Row 41
$ws.Range("E41").Value = "  +0.38%  "
# Row 42
$ws.Range("E42").Value = "  -0.17%  "
# Row 43
$ws.Range("D43").Value = "'0.982"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.22%  "
# Row 44
$ws.Range("E44").Value = "  +2.56%  "
# Row 45
$ws.Range("D45").Value = "'63.82"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.47%  "
# Row 46
$ws.Range("E46").Value = "  -1.38%  "
# Row 47
$ws.Range("D47").Value = "1.692.07"
$ws.Range("E47").Value = "  -0.28%  "
# Row 48
$ws.Range("E48").Value = "  -3.21%  "
# Row 49
$ws.Range("E49").Value = "  -0.16%  "
# Row 50
$ws.Range("E50").Value = "  +0.19%  "
# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0972"
$ws.Range("E51").Value = "  -0.79%  "
